# Apply the upstream 'cryptos list' price/volume(1h) refresh to Sheet1.
# Also swaps rows 37/38 (Monero <-> PolygonEcosystemToken reordering).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.666.19"
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").Value = "2.461.23"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "572.99"
$ws.Range("E5").Value = "  -0.61%  "
$ws.Range("D6").Value = "147.35"
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -1.70%  "
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("E10").Value = "  -0.25%  "
$ws.Range("D11").Value = "5.29"
$ws.Range("E11").Value = "  +0.08%  "
$ws.Range("E12").Value = "  -0.64%  "
$ws.Range("D13").Value = "28.99"
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("D14").Value = "0.0000176"
$ws.Range("E14").Value = "  -1.75%  "
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("D16").Value = "62.616.38"
$ws.Range("E16").Value = "  -0.43%  "
$ws.Range("D17").Value = "2.469.50"
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("E18").Value = "  -1.37%  "
$ws.Range("D19").Value = "10.89"
$ws.Range("E19").Value = "  -1.81%  "
$ws.Range("D20").Value = "325.92"
$ws.Range("E20").Value = "  -1.34%  "
$ws.Range("D22").Value = "2.17"
$ws.Range("E22").Value = "  -2.57%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").Value = "10.01"
$ws.Range("E24").Value = "  +11.46%  "
$ws.Range("D25").Value = "65.35"
$ws.Range("E25").Value = "  -1.78%  "
$ws.Range("D26").Value = "642.00"
$ws.Range("E26").Value = "  -3.68%  "
$ws.Range("D28").Value = "0.0₃0972"
$ws.Range("E28").Value = "  -3.14%  "
$ws.Range("E29").Value = "  -11.29%  "
$ws.Range("D30").Value = "1.43"
$ws.Range("E30").Value = "  -1.24%  "
$ws.Range("D31").Value = "7.93"
$ws.Range("E31").Value = "  -3.17%  "
$ws.Range("E32").Value = "  -2.38%  "
$ws.Range("E33").Value = "  -4.17%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("E35").Value = "  -1.65%  "
$ws.Range("D36").Value = "4.74"
$ws.Range("E36").Value = "  -1.16%  "
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").Value = "151.74"
$ws.Range("E37").Value = "  -0.92%  "
$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").Value = "0.368"
$ws.Range("E38").Value = "  -1.40%  "
$ws.Range("E39").Value = "  -1.64%  "
$ws.Range("D40").Value = "5.30"
$ws.Range("E40").Value = "  -3.85%  "
$ws.Range("D41").Value = "2.72"
$ws.Range("E41").Value = "  -0.73%  "
$ws.Range("D42").Value = "1.73"
$ws.Range("E42").Value = "  -2.52%  "
$ws.Range("D43").Value = "0.0₆0312"
$ws.Range("E43").Value = "  -11.09%  "
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").Value = "152.99"
$ws.Range("E45").Value = "  +4.30%  "
$ws.Range("E46").Value = "  +1.06%  "
$ws.Range("E47").Value = "  -1.53%  "
$ws.Range("D48").Value = "20.31"
$ws.Range("E48").Value = "  -2.22%  "
$ws.Range("D49").Value = "0.606"
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("E50").Value = "  -1.66%  "
$ws.Range("E51").Value = "  -1.39%  "
